$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "Wed Nov 02 16:23:34 EDT 2022"
$ws.Range("B3").Value = "Wed Nov 02 16:23:46 EDT 2022"
$ws.Range("B4").Value = "Wed Nov 02 16:23:56 EDT 2022"
$ws.Range("B5").Value = "Wed Nov 02 16:24:07 EDT 2022"
$ws.Range("B6").Value = "Wed Nov 02 16:24:17 EDT 2022"
$ws.Range("B7").Value = "Wed Nov 02 16:24:28 EDT 2022"
